# Update the "Doctors" sheet data (rows 2-6, columns A-D) with the new
# set of Chennai Gynecologist/Obstetrician doctors, replacing the former
# Mumbai/Pune/Delhi/Chennai ENT specialists list.
$wb = $excel.ActiveWorkbook
$doctors = $wb.Worksheets.Item("Doctors")
$demo = $wb.Worksheets.Item("Demo")

$doctors.Range("A2").Value = "Dr. Nazira Sadique"
$doctors.Range("B2").Value = "Gynecologist/Obstetrician"
$doctors.Range("C2").Value = "39 years experience overall"
$doctors.Range("D2").Value = "Anna Nagar,Chennai"

$doctors.Range("A3").Value = "Dr. P. V. Anuradha"
$doctors.Range("B3").Value = "Gynecologist/Obstetrician"
$doctors.Range("C3").Value = "36 years experience overall"
$doctors.Range("D3").Value = "Anna Nagar East,Chennai"

$doctors.Range("A4").Value = "Dr. Srikala Prasad"
$doctors.Range("B4").Value = "Gynecologist/Obstetrician"
$doctors.Range("C4").Value = "33 years experience overall"
$doctors.Range("D4").Value = "Vanagaram,Chennai"

$doctors.Range("A5").Value = "Dr. Parimalam Ramanathan"
$doctors.Range("B5").Value = "Gynecologist/Obstetrician"
$doctors.Range("C5").Value = "30 years experience overall"
$doctors.Range("D5").Value = "Perungudi,Chennai"

$doctors.Range("A6").Value = "Dr. M.H. Abinaya"
$doctors.Range("B6").Value = "Gynecologist/Obstetrician"
$doctors.Range("C6").Value = "22 years experience overall"
$doctors.Range("D6").Value = "T Nagar,Chennai"

# Switch the active/selected tab from "Doctors" back to "Demo".
$demo.Activate()
